# edit.ps1 — add "2022-Q4" sheet of holdings data, ahead of "2022-Q3",
# and update the "总计" (summary) sheet with a new row for 2022-Q4.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new row right under the
#    header for "2022-Q4" and push the existing quarters down by one row.
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Give the new A2 the same style as the data column below it (s="2").
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 39
$summary.Range("D2").Value = 3.9

# ----------------------------------------------------------------------
# 2) Create the new "2022-Q4" holdings sheet. Clone the "2022-Q3" sheet
#    (its neighbour) so every style (header row, column-A numbering) is
#    preserved exactly, then overwrite the cell contents in place.
# ----------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q3")
$anchor.Copy($anchor)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q4"

# The clone only has 8 data rows (rows 2-9); we need 39 (rows 2-40).
# Extend the formatting of the last existing data row down to row 40.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H40").PasteSpecial(-4122)

# Columns B:G hold text values (fund code/name and formatted numbers
# stored as strings) — force text storage so Excel doesn't silently
# coerce numeric-looking strings (e.g. "004702") into numbers.
$ws.Range("B2:G40").NumberFormat = "@"

function Set-RowData($r, $idx, $code, $name, $scale, $pos, $ratio, $value, $rank) {
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $code
    $ws.Cells.Item($r, 3).Value = $name
    $ws.Cells.Item($r, 4).Value = $scale
    $ws.Cells.Item($r, 5).Value = $pos
    $ws.Cells.Item($r, 6).Value = $ratio
    $ws.Cells.Item($r, 7).Value = $value
    $ws.Cells.Item($r, 8).Value = $rank
}

Set-RowData 2 0 "004702" "南方金融主题灵活配置混合A" "12.97" "92.71" "9.66" "1.2529" 2
Set-RowData 3 1 "010659" "民生加银质量领先混合A" "16.67" "91.93" "4.05" "0.6751" 10
Set-RowData 4 2 "013500" "南方金融主题灵活配置混合C" "4.80" "92.71" "9.66" "0.4637" 2
Set-RowData 5 3 "000408" "民生加银城镇化混合A" "4.91" "92.98" "4.14" "0.2033" 10
Set-RowData 6 4 "010795" "民生加银价值发现一年持有期混合A" "4.94" "93.32" "4.09" "0.2020" 10
Set-RowData 7 5 "011843" "民生加银内核驱动混合A" "3.94" "92.40" "4.12" "0.1623" 10
Set-RowData 8 6 "690005" "民生加银内需增长混合" "2.69" "92.98" "4.14" "0.1114" 10
Set-RowData 9 7 "210005" "金鹰主题优势混合" "3.06" "94.69" "3.56" "0.1089" 10
Set-RowData 10 8 "000082" "嘉实研究阿尔法股票" "6.50" "90.00" "1.47" "0.0956" 8
Set-RowData 11 9 "003655" "信澳新财富灵活配置混合" "3.03" "78.65" "3.05" "0.0924" 5
Set-RowData 12 10 "013676" "兴银兴慧一年持有混合A" "8.13" "23.86" "1.08" "0.0878" 6
Set-RowData 13 11 "004895" "华商鑫安灵活配置混合" "1.06" "92.66" "5.06" "0.0536" 2
Set-RowData 14 12 "009206" "兴银丰运稳益回报混合C" "3.03" "39.08" "1.74" "0.0527" 2
Set-RowData 15 13 "013677" "兴银兴慧一年持有混合C" "4.59" "23.86" "1.08" "0.0496" 6
Set-RowData 16 14 "519963" "长信利盈灵活配置混合A" "2.38" "34.67" "1.55" "0.0369" 4
Set-RowData 17 15 "009205" "兴银丰运稳益回报混合A" "1.91" "39.08" "1.74" "0.0332" 2
Set-RowData 18 16 "010660" "民生加银质量领先混合C" "0.82" "91.93" "4.05" "0.0332" 10
Set-RowData 19 17 "012245" "广发金融地产精选股票C" "0.55" "84.95" "5.80" "0.0319" 1
Set-RowData 20 18 "015453" "中欧中证500指数增强A" "1.10" "91.84" "1.73" "0.0190" 6
Set-RowData 21 19 "011105" "长信稳健均衡6个月持有期混合A" "1.66" "25.30" "1.14" "0.0189" 8
Set-RowData 22 20 "012244" "广发金融地产精选股票A" "0.31" "84.95" "5.80" "0.0180" 1
Set-RowData 23 21 "010796" "民生加银价值发现一年持有期混合C" "0.35" "93.32" "4.09" "0.0143" 10
Set-RowData 24 22 "002681" "金鹰元和灵活配置混合A" "0.30" "81.19" "4.77" "0.0143" 9
Set-RowData 25 23 "002682" "金鹰元和灵活配置混合C" "0.23" "81.19" "4.77" "0.0110" 9
Set-RowData 26 24 "001351" "诺安中证500指数增强A" "0.42" "93.99" "2.42" "0.0102" 3
Set-RowData 27 25 "002068" "东方多策略灵活配置混合C" "0.26" "87.87" "3.02" "0.0079" 4
Set-RowData 28 26 "015454" "中欧中证500指数增强C" "0.32" "91.84" "1.73" "0.0055" 6
Set-RowData 29 27 "011844" "民生加银内核驱动混合C" "0.13" "92.40" "4.12" "0.0054" 10
Set-RowData 30 28 "003186" "鹏华兴安定期开放灵活配置混合" "0.52" "20.17" "1.01" "0.0053" 10
Set-RowData 31 29 "011106" "长信稳健均衡6个月持有期混合C" "0.42" "25.30" "1.14" "0.0048" 8
Set-RowData 32 30 "510560" "国寿安保中证500ETF" "1.81" "99.21" "0.21" "0.0038" 9
Set-RowData 33 31 "005618" "融通红利机会主题精选灵活配置混合A" "0.13" "84.29" "2.93" "0.0038" 7
Set-RowData 34 32 "009706" "民生加银城镇化混合C" "0.05" "92.98" "4.14" "0.0021" 10
Set-RowData 35 33 "010355" "诺安中证500指数增强C" "0.08" "93.99" "2.42" "0.0019" 3
Set-RowData 36 34 "004791" "富荣中证500指数增强C" "0.09" "90.60" "2.00" "0.0018" 10
Set-RowData 37 35 "400023" "东方多策略灵活配置混合A" "0.03" "87.87" "3.02" "0.0009" 4
Set-RowData 38 36 "005619" "融通红利机会主题精选灵活配置混合C" "0.02" "84.29" "2.93" "0.0006" 7
Set-RowData 39 37 "004790" "富荣中证500指数增强A" "0.02" "90.60" "2.00" "0.0004" 10
Set-RowData 40 38 "519962" "长信利盈灵活配置混合C" "0.01" "34.67" "1.55" "0.0002" 4

# Drop the number-format override now that the text is committed, so the
# cells end up with the workbook's default (unstyled) formatting.
$ws.Range("B2:G40").Style = "Normal"
